# Added dummies, mode models
# Adds two dummy-variable columns (K: Dcovid, L: DGFC) to the regression
# data sheet and populates them for every existing data row (2-83).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Headers -------------------------------------------------------------
$ws.Cells.Item(1, 11).Value = "Dcovid"
$ws.Cells.Item(1, 12).Value = "DGFC"

# --- Dummy values per quarterly row (rows 2..83) --------------------------
# Dcovid = 1 for 2020Q2 - 2021Q3 (COVID period)
# DGFC   = 1 for 2008Q1 - 2010Q4 (Global Financial Crisis period)
$kVals = @(0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,1,1,1,1,1,1,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0)
$lVals = @(0,0,0,0,0,0,0,0,0,0,0,0,1,1,1,1,1,1,1,1,1,1,1,1,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0)

$firstRow = 2
$lastRow = 83

for ($r = $firstRow; $r -le $lastRow; $r++) {
    $idx = $r - $firstRow
    $ws.Cells.Item($r, 11).Value = $kVals[$idx]
    $ws.Cells.Item($r, 12).Value = $lVals[$idx]
}

# --- Formatting: General number format, centered, for the new data cells -
$dataRange = $ws.Range("K2:L83")
$dataRange.NumberFormat = "General"
$dataRange.HorizontalAlignment = -4108

# --- Selection, matching the author's final cursor position --------------
$ws.Range("M19").Select()

Write-Host "Added Dcovid/DGFC dummy columns"
